$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update matricule value from 16PK379 to 16PK389 (cell A5, merged A5:A7)
$ws.Range("A5").Value = "16PK389"

# Set new value "²" in cell A8 (merged A8:A10)
$ws.Range("A8").Value = "²"

# Update the selection to A8:A10 with active cell A8
$ws.Range("A8:A10").Select()
